$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Edil3"
$ws.Range("C2").Value = "Itgb3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.003754
$ws.Range("H2").Value = 0.011262
$ws.Range("I2").Value = 0.001049024329130219
$ws.Range("J2").Value = 0.001049024329130219
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 7.214110666666667
$ws.Range("N2").Value = 21.642332
$ws.Range("O2").Value = 0.4688823795981188
$ws.Range("P2").Value = 0.4688823795981188
$ws.Range("Q2").Value = 0.02708177144266666
$ws.Range("R2").Value = 0.243735942984
$ws.Range("S2").Value = 0.0004918690236988974
$ws.Range("T2").Value = 0.0004918690236988973

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Edil3"
$ws.Range("C3").Value = "Itgb3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.003754
$ws.Range("H3").Value = 0.011262
$ws.Range("I3").Value = 0.001049024329130219
$ws.Range("J3").Value = 0.001049024329130219
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 7.110350666666666
$ws.Range("N3").Value = 21.331052
$ws.Range("O3").Value = 0.4621384803214003
$ws.Range("P3").Value = 0.4621384803214003
$ws.Range("Q3").Value = 0.02669225640266666
$ws.Range("R3").Value = 0.240230307624
$ws.Range("S3").Value = 0.0004847945092844159
$ws.Range("T3").Value = 0.0004847945092844158

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Edil3"
$ws.Range("C4").Value = "Itgb3"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.003754
$ws.Range("H4").Value = 0.011262
$ws.Range("I4").Value = 0.001049024329130219
$ws.Range("J4").Value = 0.001049024329130219
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.061296333333333
$ws.Range("N4").Value = 3.183889
$ws.Range("O4").Value = 0.06897914008048092
$ws.Range("P4").Value = 0.06897914008048092
$ws.Range("Q4").Value = 0.003984106435333333
$ws.Range("R4").Value = 0.035856957918
$ws.Range("S4").Value = 0.00007236079614690591
$ws.Range("T4").Value = 0.0000723607961469059

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Edil3"
$ws.Range("C5").Value = "Itgb3"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2.635289666666667
$ws.Range("H5").Value = 7.905868999999999
$ws.Range("I5").Value = 0.7364099559506655
$ws.Range("J5").Value = 0.7364099559506654
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 7.214110666666667
$ws.Range("N5").Value = 21.642332
$ws.Range("O5").Value = 0.4688823795981188
$ws.Range("P5").Value = 0.4688823795981188
$ws.Range("Q5").Value = 19.01127129405644
$ws.Range("R5").Value = 171.101441646508
$ws.Range("S5").Value = 0.3452896525058939
$ws.Range("T5").Value = 0.3452896525058939

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Edil3"
$ws.Range("C6").Value = "Itgb3"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 2.635289666666667
$ws.Range("H6").Value = 7.905868999999999
$ws.Range("I6").Value = 0.7364099559506655
$ws.Range("J6").Value = 0.7364099559506654
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 7.110350666666666
$ws.Range("N6").Value = 21.331052
$ws.Range("O6").Value = 0.4621384803214003
$ws.Range("P6").Value = 0.4621384803214003
$ws.Range("Q6").Value = 18.73783363824311
$ws.Range("R6").Value = 168.640502744188
$ws.Range("S6").Value = 0.3403233779365899
$ws.Range("T6").Value = 0.3403233779365898

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Edil3"
$ws.Range("C7").Value = "Itgb3"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 2.635289666666667
$ws.Range("H7").Value = 7.905868999999999
$ws.Range("I7").Value = 0.7364099559506655
$ws.Range("J7").Value = 0.7364099559506654
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.061296333333333
$ws.Range("N7").Value = 3.183889
$ws.Range("O7").Value = 0.06897914008048092
$ws.Range("P7").Value = 0.06897914008048092
$ws.Range("Q7").Value = 2.796823260504555
$ws.Range("R7").Value = 25.17140934454099
$ws.Range("S7").Value = 0.05079692550818174
$ws.Range("T7").Value = 0.05079692550818173

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Edil3"
$ws.Range("C8").Value = "Itgb3"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.9395196666666666
$ws.Range("H8").Value = 2.818559
$ws.Range("I8").Value = 0.2625410197202043
$ws.Range("J8").Value = 0.2625410197202043
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 7.214110666666667
$ws.Range("N8").Value = 21.642332
$ws.Range("O8").Value = 0.4688823795981188
$ws.Range("P8").Value = 0.4688823795981188
$ws.Range("Q8").Value = 6.777798848843111
$ws.Range("R8").Value = 61.000189639588
$ws.Range("S8").Value = 0.1231008580685261
$ws.Range("T8").Value = 0.123100858068526

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Edil3"
$ws.Range("C9").Value = "Itgb3"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.9395196666666666
$ws.Range("H9").Value = 2.818559
$ws.Range("I9").Value = 0.2625410197202043
$ws.Range("J9").Value = 0.2625410197202043
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 7.110350666666666
$ws.Range("N9").Value = 21.331052
$ws.Range("O9").Value = 0.4621384803214003
$ws.Range("P9").Value = 0.4621384803214003
$ws.Range("Q9").Value = 6.680314288229777
$ws.Range("R9").Value = 60.122828594068
$ws.Range("S9").Value = 0.121330307875526
$ws.Range("T9").Value = 0.121330307875526

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Edil3"
$ws.Range("C10").Value = "Itgb3"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.9395196666666666
$ws.Range("H10").Value = 2.818559
$ws.Range("I10").Value = 0.2625410197202043
$ws.Range("J10").Value = 0.2625410197202043
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.061296333333333
$ws.Range("N10").Value = 3.183889
$ws.Range("O10").Value = 0.06897914008048092
$ws.Range("P10").Value = 0.06897914008048092
$ws.Range("Q10").Value = 0.9971087773278887
$ws.Range("R10").Value = 8.973978995950999
$ws.Range("S10").Value = 0.01810985377615228
$ws.Range("T10").Value = 0.01810985377615227
